# Fixing the api_key_entsoe missing value
# Applies corrected ENTSO-E border-flow figures for the affected
# timestamps (rows 26-29, 42-48) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 26-29: tiny float-precision correction on columns C and M
#     (and the System Direction column N that derives from them) ---
$ws.Range("C26").Value = 412.3750000000001
$ws.Range("M26").Value = 798.3750000000001
$ws.Range("N26").Value = -123.6750000000001

$ws.Range("C27").Value = 412.3750000000001
$ws.Range("M27").Value = 798.3750000000001
$ws.Range("N27").Value = -72.27500000000009

$ws.Range("C28").Value = 412.3750000000001
$ws.Range("M28").Value = 798.3750000000001
$ws.Range("N28").Value = -5.875000000000114

$ws.Range("C29").Value = 412.3750000000001
$ws.Range("M29").Value = 798.3750000000001
$ws.Range("N29").Value = -47.87500000000011

# --- Rows 42-45: previously-missing Deficit_RO_RS (E) data now filled in,
#     which also changes Excedent_RO_BG (B) and the totals (L, M, N) ---
$ws.Range("B42").Value = 20.55000000000001
$ws.Range("E42").Value = 206
$ws.Range("L42").Value = 970.3500000000001
$ws.Range("M42").Value = 664
$ws.Range("N42").Value = 306.3500000000001

$ws.Range("B43").Value = 20.55000000000001
$ws.Range("E43").Value = 206
$ws.Range("L43").Value = 979.55
$ws.Range("M43").Value = 664
$ws.Range("N43").Value = 315.55

$ws.Range("B44").Value = 20.55000000000001
$ws.Range("E44").Value = 206
$ws.Range("L44").Value = 960.55
$ws.Range("M44").Value = 664
$ws.Range("N44").Value = 296.55

$ws.Range("B45").Value = 20.55000000000001
$ws.Range("E45").Value = 206
$ws.Range("L45").Value = 927.55
$ws.Range("M45").Value = 664
$ws.Range("N45").Value = 263.55

# --- Rows 46-48: previously-missing Excedent_RO_HU (F) data now filled in,
#     which also changes the totals (L, N) ---
$ws.Range("F46").Value = 934
$ws.Range("L46").Value = 934
$ws.Range("N46").Value = 934

$ws.Range("F47").Value = 937
$ws.Range("L47").Value = 937
$ws.Range("N47").Value = 937

$ws.Range("F48").Value = 988
$ws.Range("L48").Value = 988
$ws.Range("N48").Value = 988

Write-Host "Applied Border_Flows corrections for rows 26-29 and 42-48"
